$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "interactionPath" column (E) with header
$ws.Range("E1").Value = "interactionPath"

# Fill in interactionPath values derived from dayID/interactionID combos
$ws.Range("E2").Value = "Cashier_Interaction_Day01_01"
$ws.Range("E3").Value = "Cashier_Interaction_Day01_02"
$ws.Range("E4").Value = "Cashier_Interaction_Day02_01"
$ws.Range("E5").Value = "Cashier_Interaction_Day02_02"

# Size the new column similarly to how it was authored in Excel
$ws.Columns.Item(5).ColumnWidth = 27.3

# Update the active selection to reflect where the author left off
$ws.Range("H11").Select()
